$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 20) with the next forecast data point,
# mirroring the style/format of the existing data rows.
$ws.Range("A20").Value = 45986
$ws.Range("A20").Style = $ws.Range("A19").Style

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.560577522109297
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.325305149734723
